$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 3
    4 = 1
    5 = 1
    6 = 2
    7 = 1
    8 = 3
    9 = 0
    10 = 3
    11 = 2
    12 = 0
    13 = 2
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 2
    25 = 4
    26 = 2
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 4
    32 = 0
    33 = 1
    34 = 2
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 4
    40 = 5
    41 = 2
    42 = 0
    43 = 2
    44 = 6
    45 = 1
    46 = 0
    47 = 1
    48 = 0
    50 = 1
    51 = 2
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 2
    57 = 2
    58 = 0
    59 = 1
    60 = 1
    61 = 2
    62 = 0
    63 = 1
    64 = 0
    65 = 2
    67 = 0
    68 = 2
    69 = 1
    70 = 2
    71 = 3
    72 = 1
    74 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
